$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22: change A22 from text "328/GFG" to number 328, and add C22 = "Java/Python"
$ws.Range("A22").Value = 328
$ws.Range("A22").HorizontalAlignment = -4131  # xlHAlignLeft, matches style used by A24/A25
$ws.Range("C22").Value = "Java/Python"

# Row 23: add C23 = "Python"
$ws.Range("C23").Value = "Python"

# Update selection to match new active cell
$ws.Range("E24").Select()
